$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix total marks error on rows 11 (Marking) and 12 (Total)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

$ws.Range("B12").Value = 40
$ws.Range("C12").Value = -6
$ws.Range("E12").Value = "34 / 112"
